# Fix wrong link to docx version of cards:
# split the run ", Gilles Babinet" into ", Gilles " and "Babinet"
# (same run formatting) on slide 2's speaker/influencer-list shape.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(14)
$tr = $shp.TextFrame.TextRange

$full = $tr.Text
$target = "Babinet"
$idx = $full.IndexOf($target)

if ($idx -ge 0) {
    $start1 = $idx + 1   # Characters() uses 1-based indexing
    $sub = $tr.Characters($start1, $target.Length)
    # Re-assigning the same text forces PowerPoint to split the run at
    # this boundary while preserving the existing run formatting (rPr).
    $sub.Text = $target
}
